$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, mirrors style of existing header row (bold/centered/bordered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# "Save" column data, rows 2-12
$saveValues = @(0, 1, 1, 0, 1, 0, 1, 0, 0, 1, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
